# Customer purchase error solved
# PV scenarios included
#
# The old pv_production sheet held a bogus 4-row "period/pv_power" lookup
# (left over from a customer-purchase placeholder). Replace it with the
# real PV-output scenario table: 16 scenarios (w1..w16) x 4 periods
# (t1..t4), matching the layout already used by the other scenario sheets
# (pool_price_scenarios, non_anticipativity_matrix, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pv_production")

# Wipe the obsolete A1:B5 table before laying down the new A1:E17 one.
$ws.Cells.Clear() | Out-Null

$ws.Range("A1").Value = "Scenario"
$ws.Range("B1").Value = "t1"
$ws.Range("C1").Value = "t2"
$ws.Range("D1").Value = "t3"
$ws.Range("E1").Value = "t4"

$data = @(
    @("w1",  44, 47, 45, 45),
    @("w2",  44, 47, 45, 47),
    @("w3",  44, 47, 46, 47),
    @("w4",  44, 47, 46, 49),
    @("w5",  44, 44, 42, 42),
    @("w6",  44, 44, 42, 45),
    @("w7",  44, 44, 44, 45),
    @("w8",  44, 44, 44, 46),
    @("w9",  50, 52, 50, 51),
    @("w10", 50, 52, 50, 55),
    @("w11", 50, 52, 51, 52),
    @("w12", 50, 52, 51, 53),
    @("w13", 50, 50, 47, 49),
    @("w14", 50, 50, 47, 50),
    @("w15", 50, 50, 48, 50),
    @("w16", 50, 50, 48, 53)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# Leave behind the same cursor positions recorded on each sheet after the edit.
$wb.Worksheets.Item("bilateral_contract_data").Range("F2").Select() | Out-Null
$wb.Worksheets.Item("pool_price_scenarios").Range("A1:E17").Select() | Out-Null
$wb.Worksheets.Item("non_anticipativity_matrix").Range("C20").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("G7").Select() | Out-Null
